$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as "379.63", "0.110", "2.037.98" (dot-
# separated thousands). These look numeric but must stay literal text,
# exactly as authored in the original inlineStr cells - e.g. "0.110" must
# not collapse to "0.11", and "2.037.98" is not even a legal number. A
# leading apostrophe is Excel's normal "force text" entry convention; it
# stops COM auto-coercing the assigned string into a Double, and the
# apostrophe itself is not stored as part of the cell's text.

$ws.Range('D2').Value = "'51.108.02"
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = "'2.958.43"
$ws.Range('E3').Value = '  +1.02%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = "'379.63"
$ws.Range('E5').Value = '  +1.74%  '
$ws.Range('D6').Value = "'102.24"
$ws.Range('E6').Value = '  +1.10%  '
$ws.Range('E7').Value = '  +1.92%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  +2.25%  '
$ws.Range('D10').Value = "'36.38"
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('E12').Value = '  +2.45%  '
$ws.Range('D13').Value = "'7.83"
$ws.Range('E13').Value = '  +6.49%  '
$ws.Range('D14').Value = "'3.423.26"
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('D15').Value = "'18.28"
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = "'2.960.06"
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').Value = "'11.11"
$ws.Range('E17').Value = '  -9.75%  '
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').Value = "'51.165.25"
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').Value = "'3.15"
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').Value = "'12.40"
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  +1.25%  '
$ws.Range('D23').Value = "'70.26"
$ws.Range('E23').Value = '  +2.89%  '
$ws.Range('D24').Value = "'266.91"
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').Value = "'3.21"
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('D26').Value = "'7.83"
$ws.Range('E26').Value = '  -1.79%  '
$ws.Range('D27').Value = "'7.31"
$ws.Range('E27').Value = '  -2.03%  '
$ws.Range('D28').Value = "'0.999"
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = "'25.88"
$ws.Range('E29').Value = '  +1.86%  '
$ws.Range('D30').Value = "'0.163"
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('D31').Value = "'0.110"
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').Value = "'10.29"
$ws.Range('E32').Value = '  +3.51%  '
$ws.Range('D33').Value = "'34.43"
$ws.Range('E33').Value = '  +5.42%  '
$ws.Range('D34').Value = "'51.13"
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').Value = "'2.06"
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('D36').Value = "'0.0434"
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').Value = "'3.25"
$ws.Range('E38').Value = '  +3.76%  '
$ws.Range('E39').Value = '  +1.33%  '
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('D41').Value = "'16.44"
$ws.Range('E41').Value = '  +2.13%  '
$ws.Range('D42').Value = "'125.17"
$ws.Range('E42').Value = '  +4.22%  '
$ws.Range('D43').Value = "'2.50"
$ws.Range('E43').Value = '  +1.95%  '
$ws.Range('D44').Value = "'3.55"
$ws.Range('E44').Value = '  +8.69%  '
$ws.Range('D45').Value = "'21.49"
$ws.Range('E45').Value = '  +3.05%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = "'2.38"
$ws.Range('E46').Value = '  +3.29%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').Value = "'0.272"
$ws.Range('E47').Value = '  -0.56%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = "'2.02"
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('D49').Value = "'2.037.98"
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('D50').Value = "'0.0322"
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('E51').Value = '  +7.01%  '
